# 15.1.2.xlsx — add a new "2020" data column (L) that duplicates the
# existing last column (K) for both the year header row (row 3) and the
# value row (row 4), then leave the selection on L10 (matches the
# author's recorded cursor position after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy K3:K4 (year header "2020" + value "6.18", with their formatting)
# into the new column L3:L4. Using Copy(Destination) duplicates both the
# values and the cell styles in one shot, so column L ends up identical
# to column K, exactly like the extra "2020" column that was appended in
# the source workbook.
$ws.Range("K3:K4").Copy($ws.Range("L3:L4"))

# Clear the marching-ants clipboard state left over from Copy().
$excel.CutCopyMode = $false

# Restore/record the selection on L10, matching the post-edit cursor.
$ws.Range("L10").Select() | Out-Null
